$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (i.e. as the 2nd tab),
#    pushing the existing "2022-Q1" / "2021-Q4" tabs one slot to the right.
#    NOTE: fetch sheet references freshly (by name) *after* Add() - this
#    runtime re-seats already-held worksheet variables by position once a
#    new sheet is inserted, so stale variables can silently point at the
#    wrong tab.
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $zongji)
$q4.Name = "2022-Q4"

$q1 = $wb.Worksheets.Item("2022-Q1")

# Clone structure + formatting of the existing "2022-Q1" sheet (same header
# layout / styles) onto the new sheet, then overwrite with the Q4 figures.
$q1.UsedRange.Copy()
$q4.Range("A1").PasteSpecial(-4104)

$q1.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

# --- 2022-Q4 fund holdings data ---
# Fund code (B) and the numeric-looking text columns (D:G) must stay text,
# exactly like the other quarter sheets (otherwise e.g. "001706" truncates
# to 1706, or "6.70" normalises to "6.7"). A leading apostrophe forces
# Excel to store the literal text, same as a user typing it in manually.

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "'001706"
$q4.Cells.Item(2,3).Value = "诺安积极回报灵活配置混合A"
$q4.Cells.Item(2,4).Value = "'0.52"
$q4.Cells.Item(2,5).Value = "'93.31"
$q4.Cells.Item(2,6).Value = "'6.70"
$q4.Cells.Item(2,7).Value = "'0.0348"
$q4.Cells.Item(2,8).Value = 9

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "'015564"
$q4.Cells.Item(3,3).Value = "大成弘远回报一年持有混合A"
$q4.Cells.Item(3,4).Value = "'2.54"
$q4.Cells.Item(3,5).Value = "'27.63"
$q4.Cells.Item(3,6).Value = "'1.22"
$q4.Cells.Item(3,7).Value = "'0.0310"
$q4.Cells.Item(3,8).Value = 8

$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "'012847"
$q4.Cells.Item(4,3).Value = "诺安积极回报灵活配置混合C"
$q4.Cells.Item(4,4).Value = "'0.18"
$q4.Cells.Item(4,5).Value = "'93.31"
$q4.Cells.Item(4,6).Value = "'6.70"
$q4.Cells.Item(4,7).Value = "'0.0121"
$q4.Cells.Item(4,8).Value = 9

$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "'015565"
$q4.Cells.Item(5,3).Value = "大成弘远回报一年持有混合C"
$q4.Cells.Item(5,4).Value = "'0.09"
$q4.Cells.Item(5,5).Value = "'27.63"
$q4.Cells.Item(5,6).Value = "'1.22"
$q4.Cells.Item(5,7).Value = "'0.0011"
$q4.Cells.Item(5,8).Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new "2022-Q4" row right under
#    the header, pushing the existing "2022-Q1" / "2021-Q4" rows down one.
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")

# Make sure the newly-used index cells (A3:A4) pick up the same style as the
# existing index column (A2) before we repopulate the rows.
$zongji.Range("A2").Copy()
$zongji.Range("A2:A4").PasteSpecial(-4122)

$zongji.Cells.Item(2,1).Value = 0
$zongji.Cells.Item(2,2).Value = "2022-Q4"
$zongji.Cells.Item(2,3).Value = 4
$zongji.Cells.Item(2,4).Value = 0.08

$zongji.Cells.Item(3,1).Value = 1
$zongji.Cells.Item(3,2).Value = "2022-Q1"
$zongji.Cells.Item(3,3).Value = 1
$zongji.Cells.Item(3,4).Value = 0.01

$zongji.Cells.Item(4,1).Value = 2
$zongji.Cells.Item(4,2).Value = "2021-Q4"
$zongji.Cells.Item(4,3).Value = 4
$zongji.Cells.Item(4,4).Value = 1.14

Write-Output "done"
